$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marks: Marking row (B11) right-answer count, Total row (B12) and
# the Corr/total marks text in E12.
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 65
$ws.Range("E12").Value = "65/140"
